# Delete the "numberOfPages" column (column E) from Sheet1, shifting the
# "documentType" column (old F) left into E, matching the bad_first_column
# fixture update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns.Item(5).Delete()

$ws.Columns.Item(5).Select()
